# CryCompanywiseStockReport_1.xlsx update
# For several groups of consecutive rows that describe the same stock item
# (identical product name in column C, and purchase price in column D),
# the per-batch figures in columns B (batch/stock code), E (sale price),
# F (quantity) and G (value) had been cyclically rotated by one position.
# This script restores/applies the correct values for each affected row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B136").Value = 63902
$ws.Range("E136").Value = 34.04
$ws.Range("F136").Value = 0
$ws.Range("G136").Value = 0

$ws.Range("B137").Value = 48654
$ws.Range("E137").Value = 38.26
$ws.Range("F137").Value = -1
$ws.Range("G137").Value = -32.02

$ws.Range("B146").Value = 64350
$ws.Range("E146").Value = 70.63
$ws.Range("F146").Value = 2
$ws.Range("G146").Value = 132.88

$ws.Range("B147").Value = 57756
$ws.Range("E147").Value = 79.37
$ws.Range("F147").Value = -100
$ws.Range("G147").Value = -6644

$ws.Range("B148").Value = 53925
$ws.Range("F148").Value = 1
$ws.Range("G148").Value = 66.44

$ws.Range("B233").Value = 48719
$ws.Range("E233").Value = 353.35
$ws.Range("F233").Value = -81
$ws.Range("G233").Value = -23955.75

$ws.Range("B234").Value = 64979
$ws.Range("E234").Value = 314.41
$ws.Range("F234").Value = 11
$ws.Range("G234").Value = 3253.25

$ws.Range("B246").Value = 64973
$ws.Range("E246").Value = 35.4
$ws.Range("F246").Value = 64
$ws.Range("G246").Value = 2131.2

$ws.Range("B247").Value = 48706
$ws.Range("E247").Value = 39.8
$ws.Range("F247").Value = -144
$ws.Range("G247").Value = -4795.2

$ws.Range("B277").Value = 61610
$ws.Range("E277").Value = 122.71
$ws.Range("F277").Value = -58
$ws.Range("G277").Value = -5957.18

$ws.Range("B278").Value = 63565
$ws.Range("E278").Value = 109.19
$ws.Range("F278").Value = 60
$ws.Range("G278").Value = 6162.6

$ws.Range("B292").Value = 63520
$ws.Range("E292").Value = 153.4
$ws.Range("F292").Value = 73
$ws.Range("G292").Value = 10532.44

$ws.Range("B293").Value = 55373
$ws.Range("E293").Value = 163.62
$ws.Range("F293").Value = -94
$ws.Range("G293").Value = -13562.32

$ws.Range("B295").Value = 63571
$ws.Range("F295").Value = 4
$ws.Range("G295").Value = 573.92

$ws.Range("B296").Value = 63531
$ws.Range("F296").Value = 80
$ws.Range("G296").Value = 11478.4

$ws.Range("B299").Value = 63510
$ws.Range("E299").Value = 50.66
$ws.Range("F299").Value = 145
$ws.Range("G299").Value = 6907.8

$ws.Range("B300").Value = 55356
$ws.Range("E300").Value = 54.04
$ws.Range("F300").Value = -158
$ws.Range("G300").Value = -7527.12

$ws.Range("B420").Value = 47097
$ws.Range("D420").Value = 112.28
$ws.Range("E420").Value = 134.16
$ws.Range("F420").Value = 15
$ws.Range("G420").Value = 1684.2

$ws.Range("B421").Value = 58047
$ws.Range("D421").Value = 105.54
$ws.Range("E421").Value = 126.1
$ws.Range("F421").Value = 42
$ws.Range("G421").Value = 4432.68

$ws.Range("B472").Value = 45695
$ws.Range("E472").Value = 23.58
$ws.Range("F472").Value = -36
$ws.Range("G472").Value = -710.28

$ws.Range("B473").Value = 64915
$ws.Range("E473").Value = 20.98
$ws.Range("F473").Value = 0
$ws.Range("G473").Value = 0

$ws.Range("B479").Value = 64927
$ws.Range("E479").Value = 17.26
$ws.Range("F479").Value = 196
$ws.Range("G479").Value = 3179.12

$ws.Range("B480").Value = 45718
$ws.Range("E480").Value = 19.38
$ws.Range("F480").Value = -294
$ws.Range("G480").Value = -4768.68

$ws.Range("B485").Value = 45709
$ws.Range("E485").Value = 15.69
$ws.Range("F485").Value = -300
$ws.Range("G485").Value = -3945

$ws.Range("B486").Value = 64925
$ws.Range("E486").Value = 13.97
$ws.Range("F486").Value = 201
$ws.Range("G486").Value = 2643.15

$ws.Range("B576").Value = 64810
$ws.Range("E576").Value = 291.22
$ws.Range("F576").Value = 6
$ws.Range("G576").Value = 1643.52

$ws.Range("B577").Value = 53319
$ws.Range("E577").Value = 310.64
$ws.Range("F577").Value = -6
$ws.Range("G577").Value = -1643.52

$ws.Range("B744").Value = 65079
$ws.Range("F744").Value = 21
$ws.Range("G744").Value = 858.27

$ws.Range("B745").Value = 65362
$ws.Range("F745").Value = 54
$ws.Range("G745").Value = 2206.98
